$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/medical-fully-insured-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")
# Fixed Value column for Extension.url row mirrors the URL (shared string), update it too
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/medical-fully-insured-indicator"
# Clear the Constraint(s) value on the root Extension row (it belongs on Extension.extension, not here)
$elements.Range("AI2").ClearContents()
